# Updated BOM to reflect 1206 chip components
#
# The BOM's passive-component "Comment" column (H) referenced 0805-package
# chip parts; the board now uses 1206-package parts, and the Schottky
# rectifier comment gains a "THD" (through-hole device) suffix. Update the
# affected comment cells accordingly and leave the selection on the last
# part that was touched (H19), matching the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value  = "100nF chip 1206"
$ws.Range("H4").Value  = "10nF chip 1206"
$ws.Range("H5").Value  = "1uF chip 1206"
$ws.Range("H6").Value  = " 2.2nF chip 1206"
$ws.Range("H7").Value  = "10k Ohms chip 1206"
$ws.Range("H8").Value  = "4.7k Ohms chip 1206"
$ws.Range("H9").Value  = "MBR120VLSFT1G Schottky Rectifier, 20V, 1A THD"
$ws.Range("H10").Value = "0 Ohm Resistor chip 1206"
$ws.Range("H19").Value = "22pF chip 1206"

$ws.Range("H19").Select()
